$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.508.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.91%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.615.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.44%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'657.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +16.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +4.63%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +7.14%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.613.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.56%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'44.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.77%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.286.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.26%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'97.155.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000261"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.02%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.611.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.34%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +9.95%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.46%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'18.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.537"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +12.09%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'513.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.65%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0000207"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +4.66%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'98.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +7.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +4.64%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.807.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.37%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.153"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.64%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.55%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'11.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.81%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.12%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +4.69%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.10%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'31.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.48%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'618.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +10.47%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'8.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +7.46%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +10.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'ImmutableX"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +10.44%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.154"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.43%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'ARBITRUM"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.934"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.46%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'USDe"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'5.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +6.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'23.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.40%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'33.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.71%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Cosmos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'8.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.88%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.02%  "
$ws.Range("E51").Style = "Normal"
